# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$updates = @{
    8  = 118
    9  = 8820
    12 = 1148
    13 = 1003
    14 = 116
    17 = 238
    18 = 277
    19 = 70
    20 = 234
    21 = 1064
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    10 = 118
    11 = 8820
    14 = 1148
    15 = 1003
    16 = 116
    19 = 238
    20 = 277
    21 = 70
    22 = 234
    23 = 1064
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
